$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most cells: plain assignment (Excel auto-detects text vs number, matching
# how the source data was produced).
# A couple of "price" cells carry a significant trailing zero (e.g. 20.10,
# 0.0510) that a normal numeric assignment would silently drop (-> 20.1 /
# 0.051). Those are entered with a leading apostrophe to force literal text,
# then ClearFormats() removes the quote-prefix formatting flag so the cell
# ends up as a plain text cell like the rest of the column.

$ws.Range("D2").Value = "27.362.07"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.655.49"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "218.13"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "0.512"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").Value = "0.0631"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").Value = "'20.10"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "1.886.90"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "1.673.32"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").Value = "67.96"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "27.360.69"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "220.48"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "6.86"
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").Value = "2.58"
$ws.Range("E22").Value = "  +5.85%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "4.45"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").Value = "9.24"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "147.43"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  +1.96%  "
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").Value = "15.89"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").Value = "'0.0510"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").Value = "3.04"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("D35").Value = "1.262.35"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "0.546"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("D39").Value = "0.843"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").Value = "0.808"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("E42").Value = "  +4.96%  "
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("D44").Value = "1.797.76"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "62.28"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").Value = "91.89"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("E48").Value = "  +25.05%  "
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").Value = "7.67"
$ws.Range("D51").Value = "0.0977"
$ws.Range("E51").Value = "  -0.25%  "
